$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "67.405.92"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.514.94"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue "D5" "589.17"
$ws.Range("E5").Value = "  -0.49%  "
Set-TextValue "D6" "170.01"
$ws.Range("E6").Value = "  -3.43%  "
$ws.Range("E7").Value = "  +0.02%  "
Set-TextValue "D8" "0.520"
$ws.Range("E8").Value = "  -2.02%  "
$ws.Range("D9").Value = "2.515.46"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -3.42%  "
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("E13").Value = "  -3.78%  "
$ws.Range("D14").Value = "3.001.63"
$ws.Range("E14").Value = "  +0.84%  "
Set-TextValue "D15" "26.02"
$ws.Range("E15").Value = "  -2.74%  "
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").Value = "67.459.60"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "2.566.13"
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("E19").Value = "  +1.10%  "
Set-TextValue "D20" "11.67"
$ws.Range("E20").Value = "  +2.17%  "
Set-TextValue "D21" "363.88"
$ws.Range("E21").Value = "  +0.56%  "
Set-TextValue "D22" "4.13"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("E23").Value = "  -3.30%  "
Set-TextValue "D24" "71.75"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("E25").Value = "  -0.05%  "
Set-TextValue "D26" "1.83"
$ws.Range("E26").Value = "  -6.38%  "
Set-TextValue "D27" "9.70"
$ws.Range("E27").Value = "  -5.27%  "
$ws.Range("D28").Value = "2.662.33"
$ws.Range("E29").Value = "  -4.70%  "
Set-TextValue "D30" "527.59"
$ws.Range("E30").Value = "  -2.73%  "
Set-TextValue "D31" "8.21"
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("E32").Value = "  -0.84%  "
Set-TextValue "D33" "1.27"
$ws.Range("E33").Value = "  -4.60%  "
Set-TextValue "D34" "0.999"
$ws.Range("E35").Value = "  -1.71%  "
Set-TextValue "D36" "157.70"
Set-TextValue "D37" "19.22"
$ws.Range("E37").Value = "  +2.36%  "
$ws.Range("E38").Value = "  -2.77%  "
Set-TextValue "D39" "18.60"
$ws.Range("E39").Value = "  -0.25%  "
Set-TextValue "D40" "1.75"
$ws.Range("E40").Value = "  -2.37%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D41" "0.341"
$ws.Range("E41").Value = "  -4.13%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D42" "5.01"
$ws.Range("E42").Value = "  -3.19%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  -3.75%  "
Set-TextValue "D45" "39.29"
$ws.Range("E45").Value = "  -1.45%  "
Set-TextValue "D46" "146.92"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("E47").Value = "  -0.85%  "
Set-TextValue "D48" "0.541"
$ws.Range("E48").Value = "  -3.43%  "
$ws.Range("D49").Value = "0.0₆0270"
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("E51").Value = "  +0.22%  "
